# Daily attendance processing - 2026-01-17 23:58:01
#
# The "Recorded By" column (G) stored the recorder list as
# "System, dnasr281@gmail.com" for a batch of session rows. This pass
# normalizes the ordering to "dnasr281@gmail.com, System" for every row
# that still has the old ordering, leaving rows that were already
# recorded differently (e.g. just "System" or just the email) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

# Recorded-by column is G; scan every used row and swap the matching ones.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
